$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.929.39'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.355.27'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.676'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.95'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  -1.68%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '59.26'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '33.70'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.86%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.29'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '2.704.45'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '16.42'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.49%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.909'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '2.356.79'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = '43.855.57'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '77.85'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '257.52'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.94'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +16.42%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.74'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.52'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.66'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '177.15'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.22'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.50'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.49'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '68.24'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +27.73%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.17'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +15.95%  '
$ws.Range('E43').Value = '  +9.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.27'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.85%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.204'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.10'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.26'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.51'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.16'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.63%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '99.76'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.53%  '
